# Fill in row 6 of the daily tracker (3rd task entry, dated 2022-01-05) and
# move the active selection from D18 to D8, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing formatting (style 2: border only) from A2 into the row-6 cells
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A6,C6,F6,G6").PasteSpecial(-4122) | Out-Null

# Copy date-format cell (style 3) from B2
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null

# Copy percentage-format cell (style 4) from E2
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null

# Copy border-only cell (style 2) into D6, then add wrap text to create the new style
$ws.Range("A2").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4122) | Out-Null
$ws.Range("D6").WrapText = $true

$excel.CutCopyMode = 0

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = (Get-Date -Year 2022 -Month 1 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C6").Value = "RPA RLOGIC"
$ws.Range("D6").Value = "1.  Integration of merging tasks of DRS, Invoice update and Saw files with Warranty  are work in progress  and it has to be done " + [char]10 + "before creating the callwise data for the P&L(already done)"
$ws.Range("E6").Value = 0.9
$ws.Range("F6").Value = "WIP"

$ws.Rows.Item(6).RowHeight = 28.8

$ws.Range("D8").Select()
